$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3450
$ws.Range("I7").Value = 900
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -6224
$ws.Range("H14").Value = 3450
$ws.Range("I14").Value = 900
$ws.Range("J14").Value = 6000
$ws.Range("K14").Value = 900
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = -709
$ws.Range("N14").Value = -6382
$ws.Range("H62").Value = 97232150
$ws.Range("I62").Value = 35724916
$ws.Range("J62").Value = 312507500
$ws.Range("K62").Value = 35724916
$ws.Range("L62").Value = 312507500
$ws.Range("M62").Value = -35724292
$ws.Range("N62").Value = -312508748
$ws.Range("H65").Value = 97232150
$ws.Range("I65").Value = 35724916
$ws.Range("J65").Value = 312507500
$ws.Range("K65").Value = 178624580
$ws.Range("L65").Value = 1562537500
$ws.Range("M65").Value = -178621460
$ws.Range("N65").Value = -1562543740
$ws.Range("H98").Value = 20928070
$ws.Range("I98").Value = 8697004
$ws.Range("J98").Value = 46502116
$ws.Range("K98").Value = 8697004
$ws.Range("L98").Value = 46502116
$ws.Range("M98").Value = -8695506
$ws.Range("N98").Value = -46505112
$ws.Range("H103").Value = 52638908
$ws.Range("I103").Value = 142857920
$ws.Range("J103").Value = 11145.833
$ws.Range("K103").Value = 428573760
$ws.Range("L103").Value = 33437.499
$ws.Range("M103").Value = -428573174
$ws.Range("N103").Value = -34609.499
$ws.Range("H122").Value = 20928070
$ws.Range("I122").Value = 8697004
$ws.Range("J122").Value = 46502116
$ws.Range("K122").Value = 26091012
$ws.Range("L122").Value = 139506348
$ws.Range("M122").Value = -26088562
$ws.Range("N122").Value = -139511248
$ws.Range("H124").Value = 54980
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 54980
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 54980
$ws.Range("N124").Value = -64800
$ws.Range("H126").Value = 45660
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 45660
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 45660
$ws.Range("N126").Value = -55540
$ws.Range("H131").Value = 6835.2856
$ws.Range("I131").Value = 1016.1667
$ws.Range("J131").Value = 11199.625
$ws.Range("K131").Value = 3048.5001
$ws.Range("L131").Value = 33598.875
$ws.Range("M131").Value = 1991.4999
$ws.Range("N131").Value = -43678.875
$ws.Range("H141").Value = 5360.625
$ws.Range("I141").Value = 3147.5
$ws.Range("J141").Value = 12000
$ws.Range("K141").Value = 9442.5
$ws.Range("L141").Value = 36000
$ws.Range("M141").Value = -4262.5
$ws.Range("N141").Value = -46360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2205.8823
$ws.Range("I63").Value = 2260
$ws.Range("J63").Value = 2128.5715
$ws.Range("K63").Value = 2260
$ws.Range("L63").Value = 2128.5715
$ws.Range("M63").Value = -1574
$ws.Range("N63").Value = -3500.5715
$ws.Range("H66").Value = 2205.8823
$ws.Range("I66").Value = 2260
$ws.Range("J66").Value = 2128.5715
$ws.Range("K66").Value = 11300
$ws.Range("L66").Value = 10642.8575
$ws.Range("M66").Value = -7868
$ws.Range("N66").Value = -17506.8575
$ws.Range("H80").Value = 36000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 36000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 36000
$ws.Range("N80").Value = -37996
$ws.Range("H83").Value = 36000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 36000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 108000
$ws.Range("N83").Value = -117984
$ws.Range("H133").Value = 34066
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 34066
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 34066
$ws.Range("N133").Value = -39126

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 466.4643
$ws.Range("I22").Value = 452.34616
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 452.34616
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -279.34616
$ws.Range("N22").Value = -996
$ws.Range("H82").Value = 18031.818
$ws.Range("I82").Value = 4140.2
$ws.Range("J82").Value = 29608.166
$ws.Range("K82").Value = 4140.2
$ws.Range("L82").Value = 29608.166
$ws.Range("M82").Value = -3757.2
$ws.Range("N82").Value = -30374.166
$ws.Range("H85").Value = 18031.818
$ws.Range("I85").Value = 4140.2
$ws.Range("J85").Value = 29608.166
$ws.Range("K85").Value = 4140.2
$ws.Range("L85").Value = 29608.166
$ws.Range("M85").Value = -2814.2
$ws.Range("N85").Value = -32260.166
$ws.Range("H86").Value = 1849.99
$ws.Range("I86").Value = 1859.5858
$ws.Range("J86").Value = 900
$ws.Range("K86").Value = 1859.5858
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = -736.5858000000001
$ws.Range("N86").Value = -3146
$ws.Range("H89").Value = 1849.99
$ws.Range("I89").Value = 1859.5858
$ws.Range("J89").Value = 900
$ws.Range("K89").Value = 9297.929
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = -3681.929
$ws.Range("N89").Value = -15732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2497.8235
$ws.Range("I62").Value = 2391.4167
$ws.Range("J62").Value = 2753.2
$ws.Range("K62").Value = 2391.4167
$ws.Range("L62").Value = 2753.2
$ws.Range("M62").Value = -1767.4167
$ws.Range("N62").Value = -4001.2
$ws.Range("H65").Value = 2497.8235
$ws.Range("I65").Value = 2391.4167
$ws.Range("J65").Value = 2753.2
$ws.Range("K65").Value = 11957.0835
$ws.Range("L65").Value = 13766
$ws.Range("M65").Value = -8837.083500000001
$ws.Range("N65").Value = -20006
$ws.Range("H68").Value = 20934.375
$ws.Range("I68").Value = 23000
$ws.Range("J68").Value = 20245.834
$ws.Range("K68").Value = 23000
$ws.Range("L68").Value = 20245.834
$ws.Range("M68").Value = -22251
$ws.Range("N68").Value = -21743.834
$ws.Range("H71").Value = 20934.375
$ws.Range("I71").Value = 23000
$ws.Range("J71").Value = 20245.834
$ws.Range("K71").Value = 69000
$ws.Range("L71").Value = 60737.50199999999
$ws.Range("M71").Value = -65256
$ws.Range("N71").Value = -68225.50199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2705.4285
$ws.Range("I75").Value = 2684.5
$ws.Range("J75").Value = 2733.3333
$ws.Range("K75").Value = 8053.5
$ws.Range("L75").Value = 8199.999899999999
$ws.Range("M75").Value = -7055.5
$ws.Range("N75").Value = -10195.9999
$ws.Range("H78").Value = 2705.4285
$ws.Range("I78").Value = 2684.5
$ws.Range("J78").Value = 2733.3333
$ws.Range("K78").Value = 24160.5
$ws.Range("L78").Value = 24599.9997
$ws.Range("M78").Value = -19168.5
$ws.Range("N78").Value = -34583.9997
$ws.Range("H131").Value = 29558.285
$ws.Range("I131").Value = 583.3333
$ws.Range("J131").Value = 32274.688
$ws.Range("K131").Value = 1749.9999
$ws.Range("L131").Value = 96824.064
$ws.Range("M131").Value = 3290.0001
$ws.Range("N131").Value = -106904.064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11675606
$ws.Range("I132").Value = 20636124
$ws.Range("J132").Value = 5350534.5
$ws.Range("K132").Value = 61908372
$ws.Range("L132").Value = 16051603.5
$ws.Range("M132").Value = -61905842
$ws.Range("N132").Value = -16056663.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2775.139
$ws.Range("I68").Value = 2681
$ws.Range("J68").Value = 2989.0908
$ws.Range("K68").Value = 2681
$ws.Range("L68").Value = 2989.0908
$ws.Range("M68").Value = -1932
$ws.Range("N68").Value = -4487.0908
$ws.Range("H71").Value = 2775.139
$ws.Range("I71").Value = 2681
$ws.Range("J71").Value = 2989.0908
$ws.Range("K71").Value = 13405
$ws.Range("L71").Value = 14945.454
$ws.Range("M71").Value = -9661
$ws.Range("N71").Value = -22433.454
$ws.Range("H93").Value = 28632
$ws.Range("I93").Value = 32951.5
$ws.Range("J93").Value = 26712.223
$ws.Range("K93").Value = 32951.5
$ws.Range("L93").Value = 26712.223
$ws.Range("M93").Value = -31703.5
$ws.Range("N93").Value = -29208.223
$ws.Range("H100").Value = 1744.6666
$ws.Range("I100").Value = 1466.6666
$ws.Range("J100").Value = 2022.6666
$ws.Range("K100").Value = 1466.6666
$ws.Range("L100").Value = 2022.6666
$ws.Range("M100").Value = -925.6666
$ws.Range("N100").Value = -3104.6666
$ws.Range("H133").Value = 63333.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 63333.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 63333.332
$ws.Range("N133").Value = -68393.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 35625
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 40285.715
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 80571.42999999999
$ws.Range("M81").Value = -4939
$ws.Range("N81").Value = -82693.42999999999
$ws.Range("H84").Value = 35625
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 40285.715
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 402857.15
$ws.Range("M84").Value = -24696
$ws.Range("N84").Value = -413465.15
$ws.Range("H103").Value = 28801
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 28801
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 28801
$ws.Range("N103").Value = -31145
$ws.Range("H133").Value = 52104.6
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 52104.6
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 52104.6
$ws.Range("N133").Value = -62224.6

Write-Output "applied changes"